# Weekly update: a new week's price report row is inserted at row 7,
# pushing the existing historical rows (7-73) down by one (8-74).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 7 (shifts rows 7:73 down to 8:74,
# and grows the sheet dimension to A1:R74 automatically).
$ws.Rows(7).Insert()

# Populate the newly inserted row 7 with this week's data.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C7").Value = 'Arica y Parinacota'
$ws.Range("D7").Value2 = 44616
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112038
$ws.Range("G7").Value = 'Cebollín baby'
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = 3250
$ws.Range("N7").Value = '$/paquete 1,5 a 2 kilos'
$ws.Range("O7").Value = 'Región de Arica y Parinacota'
$ws.Range("P7").Value = 1625
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 'Hortaliza'
